$wb = $excel.ActiveWorkbook
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# hunk 0: sheet ALC row 33
$ws_ALC.Cells.Item(33, 9).Value = 50170.4  # I33: 53750.855 -> 50170.4
$ws_ALC.Cells.Item(33, 10).Value = 605.75  # J33: 493.4 -> 605.75
$ws_ALC.Cells.Item(33, 11).Value = 50170.4  # K33: 53750.855 -> 50170.4
$ws_ALC.Cells.Item(33, 12).Value = 605.75  # L33: 493.4 -> 605.75
$ws_ALC.Cells.Item(33, 13).Value = -49941.4  # M33: -53521.855 -> -49941.4
$ws_ALC.Cells.Item(33, 14).Value = -1063.75  # N33: -951.4 -> -1063.75

# hunk 1: sheet ALC row 53
$ws_ALC.Cells.Item(53, 8).Value = 43722.434  # H53: 37263.48 -> 43722.434
$ws_ALC.Cells.Item(53, 9).Value = 210.14285  # I53: 188 -> 210.14285
$ws_ALC.Cells.Item(53, 10).Value = 62759.062  # J53: 52874.21 -> 62759.062
$ws_ALC.Cells.Item(53, 11).Value = 210.14285  # K53: 188 -> 210.14285
$ws_ALC.Cells.Item(53, 12).Value = 62759.062  # L53: 52874.21 -> 62759.062
$ws_ALC.Cells.Item(53, 13).Value = 426.85715  # M53: 449 -> 426.85715
$ws_ALC.Cells.Item(53, 14).Value = -64033.062  # N53: -54148.21 -> -64033.062

# hunk 2: sheet ALC row 87
$ws_ALC.Cells.Item(87, 8).Value = 4510000  # H87: 2282777.5 -> 4510000
$ws_ALC.Cells.Item(87, 10).Value = 9000000  # J87: 3037036.8 -> 9000000
$ws_ALC.Cells.Item(87, 12).Value = 9000000  # L87: 3037036.8 -> 9000000
$ws_ALC.Cells.Item(87, 14).Value = -9002496  # N87: -3039532.8 -> -9002496

# hunk 3: sheet ALC row 90
$ws_ALC.Cells.Item(90, 8).Value = 4510000  # H90: 2282777.5 -> 4510000
$ws_ALC.Cells.Item(90, 10).Value = 9000000  # J90: 3037036.8 -> 9000000
$ws_ALC.Cells.Item(90, 12).Value = 27000000  # L90: 9111110.399999999 -> 27000000
$ws_ALC.Cells.Item(90, 14).Value = -27012480  # N90: -9123590.399999999 -> -27012480

# hunk 4: sheet ALC row 100
$ws_ALC.Cells.Item(100, 8).Value = 8720.4  # H100: 8109.4546 -> 8720.4
$ws_ALC.Cells.Item(100, 9).Value = 12440.8  # I100: 10700.667 -> 12440.8
$ws_ALC.Cells.Item(100, 11).Value = 12440.8  # K100: 10700.667 -> 12440.8
$ws_ALC.Cells.Item(100, 13).Value = -11899.8  # M100: -10159.667 -> -11899.8

# hunk 5: sheet ALC row 113
$ws_ALC.Cells.Item(113, 8).Value = 100003000  # H113: 71431570 -> 100003000
$ws_ALC.Cells.Item(113, 10).Value = 3713.2856  # J113: 3454.4546 -> 3713.2856
$ws_ALC.Cells.Item(113, 12).Value = 3713.2856  # L113: 3454.4546 -> 3713.2856
$ws_ALC.Cells.Item(113, 14).Value = -10221.2856  # N113: -9962.454600000001 -> -10221.2856

# hunk 6: sheet ALC row 116
$ws_ALC.Cells.Item(116, 8).Value = 19486678  # H116: 17194334 -> 19486678
$ws_ALC.Cells.Item(116, 9).Value = 11962185  # I116: 10922126 -> 11962185
$ws_ALC.Cells.Item(116, 10).Value = 37043824  # J116: 30308948 -> 37043824
$ws_ALC.Cells.Item(116, 11).Value = 11962185  # K116: 10922126 -> 11962185
$ws_ALC.Cells.Item(116, 12).Value = 37043824  # L116: 30308948 -> 37043824
$ws_ALC.Cells.Item(116, 13).Value = -11958743  # M116: -10918684 -> -11958743
$ws_ALC.Cells.Item(116, 14).Value = -37050708  # N116: -30315832 -> -37050708

# hunk 7: sheet ARM row 2
$ws_ARM.Cells.Item(2, 8).Value = 618.4761999999999  # H2: 594.9091 -> 618.4761999999999
$ws_ARM.Cells.Item(2, 9).Value = 619.4  # I2: 594.6667 -> 619.4
$ws_ARM.Cells.Item(2, 11).Value = 619.4  # K2: 594.6667 -> 619.4
$ws_ARM.Cells.Item(2, 13).Value = -506.4  # M2: -481.6667 -> -506.4

# hunk 8: sheet ARM row 45
$ws_ARM.Cells.Item(45, 8).Value = 1394  # H45: 1369.3182 -> 1394
$ws_ARM.Cells.Item(45, 9).Value = 1344.4445  # I45: 1318.4736 -> 1344.4445
$ws_ARM.Cells.Item(45, 11).Value = 1344.4445  # K45: 1318.4736 -> 1344.4445
$ws_ARM.Cells.Item(45, 13).Value = -967.4445000000001  # M45: -941.4736 -> -967.4445000000001

# hunk 9: sheet ARM row 63
$ws_ARM.Cells.Item(63, 8).Value = 92313576  # H63: 109096680 -> 92313576
$ws_ARM.Cells.Item(63, 9).Value = 333335500  # I63: 500001760 -> 333335500
$ws_ARM.Cells.Item(63, 10).Value = 20007000  # J63: 22228888 -> 20007000
$ws_ARM.Cells.Item(63, 11).Value = 333335500  # K63: 500001760 -> 333335500
$ws_ARM.Cells.Item(63, 12).Value = 20007000  # L63: 22228888 -> 20007000
$ws_ARM.Cells.Item(63, 13).Value = -333334814  # M63: -500001074 -> -333334814
$ws_ARM.Cells.Item(63, 14).Value = -20008372  # N63: -22230260 -> -20008372

# hunk 10: sheet ARM row 66
$ws_ARM.Cells.Item(66, 8).Value = 92313576  # H66: 109096680 -> 92313576
$ws_ARM.Cells.Item(66, 9).Value = 333335500  # I66: 500001760 -> 333335500
$ws_ARM.Cells.Item(66, 10).Value = 20007000  # J66: 22228888 -> 20007000
$ws_ARM.Cells.Item(66, 11).Value = 1666677500  # K66: 2500008800 -> 1666677500
$ws_ARM.Cells.Item(66, 12).Value = 100035000  # L66: 111144440 -> 100035000
$ws_ARM.Cells.Item(66, 13).Value = -1666674068  # M66: -2500005368 -> -1666674068
$ws_ARM.Cells.Item(66, 14).Value = -100041864  # N66: -111151304 -> -100041864

# hunk 11: sheet ARM row 110
$ws_ARM.Cells.Item(110, 8).Value = 50002156  # H110: 47621150 -> 50002156
$ws_ARM.Cells.Item(110, 9).Value = 71430536  # I110: 71430560 -> 71430536
$ws_ARM.Cells.Item(110, 10).Value = 2604.6667  # J110: 2332.4285 -> 2604.6667
$ws_ARM.Cells.Item(110, 11).Value = 71430536  # K110: 71430560 -> 71430536
$ws_ARM.Cells.Item(110, 12).Value = 2604.6667  # L110: 2332.4285 -> 2604.6667
$ws_ARM.Cells.Item(110, 13).Value = -71428491  # M110: -71428515 -> -71428491
$ws_ARM.Cells.Item(110, 14).Value = -6694.6667  # N110: -6422.4285 -> -6694.6667

# hunk 12: sheet ARM row 116
$ws_ARM.Cells.Item(116, 8).Value = 618.4761999999999  # H116: 594.9091 -> 618.4761999999999
$ws_ARM.Cells.Item(116, 9).Value = 619.4  # I116: 594.6667 -> 619.4
$ws_ARM.Cells.Item(116, 11).Value = 619.4  # K116: 594.6667 -> 619.4
$ws_ARM.Cells.Item(116, 13).Value = 1674.6  # M116: 1699.3333 -> 1674.6

# hunk 13: sheet ARM row 122
$ws_ARM.Cells.Item(122, 8).Value = 3420.3242  # H122: 3903.25 -> 3420.3242
$ws_ARM.Cells.Item(122, 9).Value = 3019.4092  # I122: 3651.647 -> 3019.4092
$ws_ARM.Cells.Item(122, 10).Value = 4008.3333  # J122: 4188.4 -> 4008.3333
$ws_ARM.Cells.Item(122, 11).Value = 9058.2276  # K122: 10954.941 -> 9058.2276
$ws_ARM.Cells.Item(122, 12).Value = 12024.9999  # L122: 12565.2 -> 12024.9999
$ws_ARM.Cells.Item(122, 13).Value = -6608.2276  # M122: -8504.940999999999 -> -6608.2276
$ws_ARM.Cells.Item(122, 14).Value = -16924.9999  # N122: -17465.2 -> -16924.9999

# hunk 14: sheet ARM row 131
$ws_ARM.Cells.Item(131, 8).Value = 0  # H131: 125000 -> 0
$ws_ARM.Cells.Item(131, 10).Value = 0  # J131: 125000 -> 0
$ws_ARM.Cells.Item(131, 12).Value = 0  # L131: 125000 -> 0
$ws_ARM.Cells.Item(131, 14).ClearContents()  # N131: -135080 -> (removed)

# hunk 15: sheet BSM row 3
$ws_BSM.Cells.Item(3, 8).Value = 618.4761999999999  # H3: 594.9091 -> 618.4761999999999
$ws_BSM.Cells.Item(3, 9).Value = 619.4  # I3: 594.6667 -> 619.4
$ws_BSM.Cells.Item(3, 11).Value = 619.4  # K3: 594.6667 -> 619.4
$ws_BSM.Cells.Item(3, 13).Value = -505.4  # M3: -480.6667 -> -505.4

# hunk 16: sheet BSM row 105
$ws_BSM.Cells.Item(105, 8).Value = 2087.946  # H105: 2128.7715 -> 2087.946
$ws_BSM.Cells.Item(105, 9).Value = 1879.9615  # I105: 1922.1666 -> 1879.9615
$ws_BSM.Cells.Item(105, 11).Value = 1879.9615  # K105: 1922.1666 -> 1879.9615
$ws_BSM.Cells.Item(105, 13).Value = -132.9614999999999  # M105: -175.1666 -> -132.9614999999999

# hunk 17: sheet CRP row 99
$ws_CRP.Cells.Item(99, 8).Value = 6499.9  # H99: 7437.5 -> 6499.9
$ws_CRP.Cells.Item(99, 9).Value = 7999.8335  # I99: 10625 -> 7999.8335
$ws_CRP.Cells.Item(99, 11).Value = 7999.8335  # K99: 10625 -> 7999.8335
$ws_CRP.Cells.Item(99, 13).Value = -6501.8335  # M99: -9127 -> -6501.8335

# hunk 18: sheet CRP row 116
$ws_CRP.Cells.Item(116, 8).Value = 34444  # H116: 0 -> 34444
$ws_CRP.Cells.Item(116, 10).Value = 34444  # J116: 0 -> 34444
$ws_CRP.Cells.Item(116, 12).Value = 34444  # L116: 0 -> 34444
$ws_CRP.Cells.Item(116, 14).Value = -43622  # N116: None -> -43622

# hunk 19: sheet CRP row 126
$ws_CRP.Cells.Item(126, 8).Value = 6499.9  # H126: 7437.5 -> 6499.9
$ws_CRP.Cells.Item(126, 9).Value = 7999.8335  # I126: 10625 -> 7999.8335
$ws_CRP.Cells.Item(126, 11).Value = 23999.5005  # K126: 31875 -> 23999.5005
$ws_CRP.Cells.Item(126, 13).Value = -21529.5005  # M126: -29405 -> -21529.5005

# hunk 20: sheet CRP row 132
$ws_CRP.Cells.Item(132, 8).Value = 529364.8  # H132: 502966.6 -> 529364.8
$ws_CRP.Cells.Item(132, 9).Value = 2362.3333  # I132: 2302.1875 -> 2362.3333
$ws_CRP.Cells.Item(132, 11).Value = 7086.999899999999  # K132: 6906.5625 -> 7086.999899999999
$ws_CRP.Cells.Item(132, 13).Value = -4556.999899999999  # M132: -4376.5625 -> -4556.999899999999

# hunk 21: sheet CUL row 81
$ws_CUL.Cells.Item(81, 8).Value = 9838.462  # H81: 8734.916999999999 -> 9838.462
$ws_CUL.Cells.Item(81, 10).Value = 11433.333  # J81: 9977.375 -> 11433.333
$ws_CUL.Cells.Item(81, 12).Value = 34299.999  # L81: 29932.125 -> 34299.999
$ws_CUL.Cells.Item(81, 14).Value = -36545.999  # N81: -32178.125 -> -36545.999

# hunk 22: sheet CUL row 84
$ws_CUL.Cells.Item(84, 8).Value = 9838.462  # H84: 8734.916999999999 -> 9838.462
$ws_CUL.Cells.Item(84, 10).Value = 11433.333  # J84: 9977.375 -> 11433.333
$ws_CUL.Cells.Item(84, 12).Value = 102899.997  # L84: 89796.375 -> 102899.997
$ws_CUL.Cells.Item(84, 14).Value = -114131.997  # N84: -101028.375 -> -114131.997

# hunk 23: sheet CUL row 132
$ws_CUL.Cells.Item(132, 8).Value = 2592.5212  # H132: 2618.257 -> 2592.5212
$ws_CUL.Cells.Item(132, 9).Value = 900  # I132: 933.3333 -> 900
$ws_CUL.Cells.Item(132, 10).Value = 2693.5671  # J132: 2693.7014 -> 2693.5671
$ws_CUL.Cells.Item(132, 11).Value = 8100  # K132: 8399.9997 -> 8100
$ws_CUL.Cells.Item(132, 12).Value = 24242.1039  # L132: 24243.3126 -> 24242.1039
$ws_CUL.Cells.Item(132, 13).Value = -5570  # M132: -5869.9997 -> -5570
$ws_CUL.Cells.Item(132, 14).Value = -29302.1039  # N132: -29303.3126 -> -29302.1039

# hunk 24: sheet GSM row 10
$ws_GSM.Cells.Item(10, 8).Value = 9999  # H10: 7499.5 -> 9999
$ws_GSM.Cells.Item(10, 9).Value = 9999  # I10: 7499.5 -> 9999
$ws_GSM.Cells.Item(10, 11).Value = 9999  # K10: 7499.5 -> 9999
$ws_GSM.Cells.Item(10, 13).Value = -9830  # M10: -7330.5 -> -9830

# hunk 25: sheet GSM row 11
$ws_GSM.Cells.Item(11, 8).Value = 8496158  # H11: 9927857 -> 8496158
$ws_GSM.Cells.Item(11, 9).Value = 3485210  # I11: 0 -> 3485210
$ws_GSM.Cells.Item(11, 11).Value = 3485210  # K11: 0 -> 3485210
$ws_GSM.Cells.Item(11, 13).Value = -3485071  # M11: None -> -3485071

# hunk 26: sheet GSM row 12
$ws_GSM.Cells.Item(12, 8).Value = 1999.8235  # H12: 1999.9412 -> 1999.8235
$ws_GSM.Cells.Item(12, 9).Value = 2000  # I12: 1999.5 -> 2000
$ws_GSM.Cells.Item(12, 10).Value = 1999.8125  # J12: 2000 -> 1999.8125
$ws_GSM.Cells.Item(12, 11).Value = 2000  # K12: 1999.5 -> 2000
$ws_GSM.Cells.Item(12, 12).Value = 1999.8125  # L12: 2000 -> 1999.8125
$ws_GSM.Cells.Item(12, 14).Value = -2279.8125  # N12: -2280 -> -2279.8125

# hunk 27: sheet GSM row 97
$ws_GSM.Cells.Item(97, 8).Value = 614.6667  # H97: 581.8333 -> 614.6667
$ws_GSM.Cells.Item(97, 9).Value = 505.66666  # I97: 496.6 -> 505.66666
$ws_GSM.Cells.Item(97, 10).Value = 832.6667  # J97: 752.3 -> 832.6667
$ws_GSM.Cells.Item(97, 11).Value = 505.66666  # K97: 496.6 -> 505.66666
$ws_GSM.Cells.Item(97, 12).Value = 832.6667  # L97: 752.3 -> 832.6667
$ws_GSM.Cells.Item(97, 13).Value = -9.666659999999979  # M97: -0.6000000000000227 -> -9.666659999999979
$ws_GSM.Cells.Item(97, 14).Value = -1824.6667  # N97: -1744.3 -> -1824.6667

# hunk 28: sheet GSM row 122
$ws_GSM.Cells.Item(122, 8).Value = 3142.7144  # H122: 1881.7391 -> 3142.7144
$ws_GSM.Cells.Item(122, 9).Value = 2874.75  # I122: 1777.6428 -> 2874.75
$ws_GSM.Cells.Item(122, 10).Value = 3500  # J122: 2043.6666 -> 3500
$ws_GSM.Cells.Item(122, 11).Value = 8624.25  # K122: 5332.928400000001 -> 8624.25
$ws_GSM.Cells.Item(122, 12).Value = 10500  # L122: 6130.9998 -> 10500
$ws_GSM.Cells.Item(122, 13).Value = -6174.25  # M122: -2882.928400000001 -> -6174.25
$ws_GSM.Cells.Item(122, 14).Value = -15400  # N122: -11030.9998 -> -15400

# hunk 29: sheet GSM row 132
$ws_GSM.Cells.Item(132, 8).Value = 2547.4443  # H132: 2626.3845 -> 2547.4443
$ws_GSM.Cells.Item(132, 9).Value = 2331.36  # I132: 2407.875 -> 2331.36
$ws_GSM.Cells.Item(132, 11).Value = 6994.08  # K132: 7223.625 -> 6994.08
$ws_GSM.Cells.Item(132, 13).Value = -4464.08  # M132: -4693.625 -> -4464.08

# hunk 30: sheet LTW row 43
$ws_LTW.Cells.Item(43, 8).Value = 5338800  # H43: 4713470.5 -> 5338800
$ws_LTW.Cells.Item(43, 9).Value = 3131500  # I43: 2785444.5 -> 3131500
$ws_LTW.Cells.Item(43, 10).Value = 7861428.5  # J43: 6882500 -> 7861428.5
$ws_LTW.Cells.Item(43, 11).Value = 3131500  # K43: 2785444.5 -> 3131500
$ws_LTW.Cells.Item(43, 12).Value = 7861428.5  # L43: 6882500 -> 7861428.5
$ws_LTW.Cells.Item(43, 13).Value = -3131307  # M43: -2785251.5 -> -3131307
$ws_LTW.Cells.Item(43, 14).Value = -7861814.5  # N43: -6882886 -> -7861814.5

# hunk 31: sheet LTW row 46
$ws_LTW.Cells.Item(46, 8).Value = 2471.8333  # H46: 2366.2 -> 2471.8333
$ws_LTW.Cells.Item(46, 10).Value = 3210  # J46: 3280 -> 3210
$ws_LTW.Cells.Item(46, 12).Value = 3210  # L46: 3280 -> 3210
$ws_LTW.Cells.Item(46, 14).Value = -3586  # N46: -3656 -> -3586

# hunk 32: sheet LTW row 68
$ws_LTW.Cells.Item(68, 8).Value = 4051.7222  # H68: 5334.5454 -> 4051.7222
$ws_LTW.Cells.Item(68, 9).Value = 2848  # I68: 3536 -> 2848
$ws_LTW.Cells.Item(68, 10).Value = 5255.4443  # J68: 6833.3335 -> 5255.4443
$ws_LTW.Cells.Item(68, 11).Value = 2848  # K68: 3536 -> 2848
$ws_LTW.Cells.Item(68, 12).Value = 5255.4443  # L68: 6833.3335 -> 5255.4443
$ws_LTW.Cells.Item(68, 13).Value = -2099  # M68: -2787 -> -2099
$ws_LTW.Cells.Item(68, 14).Value = -6753.4443  # N68: -8331.333500000001 -> -6753.4443

# hunk 33: sheet LTW row 69
$ws_LTW.Cells.Item(69, 8).Value = 59999  # H69: 0 -> 59999
$ws_LTW.Cells.Item(69, 9).Value = 59999  # I69: 0 -> 59999
$ws_LTW.Cells.Item(69, 11).Value = 59999  # K69: 0 -> 59999
$ws_LTW.Cells.Item(69, 13).Value = -59188  # M69: None -> -59188

# hunk 34: sheet LTW row 71
$ws_LTW.Cells.Item(71, 8).Value = 4051.7222  # H71: 5334.5454 -> 4051.7222
$ws_LTW.Cells.Item(71, 9).Value = 2848  # I71: 3536 -> 2848
$ws_LTW.Cells.Item(71, 10).Value = 5255.4443  # J71: 6833.3335 -> 5255.4443
$ws_LTW.Cells.Item(71, 11).Value = 14240  # K71: 17680 -> 14240
$ws_LTW.Cells.Item(71, 12).Value = 26277.2215  # L71: 34166.6675 -> 26277.2215
$ws_LTW.Cells.Item(71, 13).Value = -10496  # M71: -13936 -> -10496
$ws_LTW.Cells.Item(71, 14).Value = -33765.2215  # N71: -41654.6675 -> -33765.2215

# hunk 35: sheet LTW row 72
$ws_LTW.Cells.Item(72, 8).Value = 59999  # H72: 0 -> 59999
$ws_LTW.Cells.Item(72, 9).Value = 59999  # I72: 0 -> 59999
$ws_LTW.Cells.Item(72, 11).Value = 179997  # K72: 0 -> 179997
$ws_LTW.Cells.Item(72, 12).Value = 0  # L72: None -> 0
$ws_LTW.Cells.Item(72, 13).Value = -175941  # M72: None -> -175941

# hunk 36: sheet WVR row 15
$ws_WVR.Cells.Item(15, 8).Value = 10000  # H15: 9000 -> 10000
$ws_WVR.Cells.Item(15, 9).Value = 0  # I15: 6000 -> 0
$ws_WVR.Cells.Item(15, 11).Value = 0  # K15: 6000 -> 0
$ws_WVR.Cells.Item(15, 13).ClearContents()  # M15: -5712 -> (removed)

# hunk 37: sheet WVR row 31
$ws_WVR.Cells.Item(31, 8).Value = 13499.5  # H31: 25999 -> 13499.5
$ws_WVR.Cells.Item(31, 9).Value = 1000  # I31: 0 -> 1000
$ws_WVR.Cells.Item(31, 11).Value = 1000  # K31: 0 -> 1000
$ws_WVR.Cells.Item(31, 13).Value = -652  # M31: None -> -652

# hunk 38: sheet WVR row 96
$ws_WVR.Cells.Item(96, 8).Value = 2369.8  # H96: 3058.1667 -> 2369.8
$ws_WVR.Cells.Item(96, 10).Value = 3316.3333  # J96: 4112.25 -> 3316.3333
$ws_WVR.Cells.Item(96, 12).Value = 3316.3333  # L96: 4112.25 -> 3316.3333
$ws_WVR.Cells.Item(96, 14).Value = -6062.3333  # N96: -6858.25 -> -6062.3333

# hunk 39: sheet WVR row 117
$ws_WVR.Cells.Item(117, 8).Value = 32989.668  # H117: 52249.75 -> 32989.668
$ws_WVR.Cells.Item(117, 10).Value = 32989.668  # J117: 52249.75 -> 32989.668
$ws_WVR.Cells.Item(117, 12).Value = 32989.668  # L117: 52249.75 -> 32989.668
$ws_WVR.Cells.Item(117, 14).Value = -42167.668  # N117: -61427.75 -> -42167.668

# hunk 40: sheet WVR row 130
$ws_WVR.Cells.Item(130, 8).Value = 51214.5  # H130: 52429 -> 51214.5
$ws_WVR.Cells.Item(130, 10).Value = 51214.5  # J130: 52429 -> 51214.5
$ws_WVR.Cells.Item(130, 12).Value = 51214.5  # L130: 52429 -> 51214.5
$ws_WVR.Cells.Item(130, 14).Value = -61254.5  # N130: -62469 -> -61254.5

# hunk 41: sheet WVR row 132
$ws_WVR.Cells.Item(132, 8).Value = 332439.6  # H132: 327089.6 -> 332439.6
$ws_WVR.Cells.Item(132, 9).Value = 438919.84  # I132: 429596.84 -> 438919.84
$ws_WVR.Cells.Item(132, 11).Value = 1316759.52  # K132: 1288790.52 -> 1316759.52
$ws_WVR.Cells.Item(132, 13).Value = -1314229.52  # M132: -1286260.52 -> -1314229.52

# hunk 42: sheet WVR row 136
$ws_WVR.Cells.Item(136, 8).Value = 4977.543  # H136: 5103.5586 -> 4977.543
$ws_WVR.Cells.Item(136, 9).Value = 5089.893  # I136: 5252.7407 -> 5089.893
$ws_WVR.Cells.Item(136, 11).Value = 15269.679  # K136: 15758.2221 -> 15269.679
$ws_WVR.Cells.Item(136, 13).Value = -12719.679  # M136: -13208.2221 -> -12719.679
